$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting the existing rows 62:154 down to 63:155.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record
# (date 2021-12-10 / Calameño / Primera, Vega Monumental Concepción - Melón).
$ws.Range("A62").Value = 11
$ws.Range("B62").Value = "Vega Monumental Concepción"
$ws.Range("C62").Value = "Bíobío"
$ws.Range("D62").Value = 44540
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = 100112027
$ws.Range("G62").Value = "Melón"
$ws.Range("H62").Value = "Calameño"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 1100
$ws.Range("M62").Value = 1044
$ws.Range("N62").Value = "`$/unidad"
$ws.Range("O62").Value = "Región de O'Higgins"
$ws.Range("P62").Value = 1044
$ws.Range("Q62").Value = 1
$ws.Range("R62").Value = "Hortaliza"
